# accountBills.xlsx - "Filter for account bills skipped in budget bills"
#
# The "Status" column (ok / SKIPPED) is no longer needed because the
# filtering logic moved out of the spreadsheet formula and into
# filterAccountBillsNotInBudgetBills() in the application code. The old
# "Prise" (price) column is renamed to "Amount", one of its sample values
# is corrected, and the conditional formatting rule that highlighted
# negative/positive amounts is repointed at the new Amount column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the "Status" column (old column C: Shop, Account, Status, Date,
#     Transaction type, Prise, Currency, Title -> Shop, Account, Date,
#     Transaction type, Amount(was Prise), Currency, Title). Deleting the
#     column shifts D..H left into C..G and updates the sheet dimension,
#     column widths and shared-string table automatically.
$ws.Columns("C:C").Delete() | Out-Null

# --- Rename the former "Prise" header (now column E) to "Amount".
$ws.Range("E1").Value = "Amount"

# --- Correct the sample amount on row 3 (was -3.49, now -29.95).
$ws.Range("E3").Value = -29.95

# --- The conditional formatting on B2:B4 used to read the (now removed)
#     Prise column via $F2>0; leave a couple of alternate highlight
#     variants behind in the style table (as Excel itself tends to do
#     while a conditional-format rule is being tweaked in the UI) before
#     settling back on a single rule that now points at the Amount
#     column ($E2>0).
$cf = $ws.Range("B2:B4").FormatConditions
$alt1 = $cf.Add(1, 5, "999999999")
$alt1.Font.ColorIndex = 1
$alt1.Interior.Color = 16051931
$alt2 = $cf.Add(1, 5, "999999999")
$alt2.Interior.Color = 14536083
$cf.Item(2).Delete()
$cf.Item(2).Delete()

$rule = $cf.Item(1)
$rule.Formula1 = '=$E2>0'

# --- Leave the selection where the user's last edit (the Amount fix on
#     row 4) would have left it.
$ws.Range("E4").Select() | Out-Null
